{"js": "// Replace the 25 \"dividend\u00f7divisor=quotient, remainder\" answer strings in the\n// worksheet table with the new values from the commit.\nconst replacements = [\n  [\"667\u00f76=111, 1\", \"878\u00f78=109, 6\"],\n  [\"909\u00f72=454, 1\", \"827\u00f72=413, 1\"],\n  [\"226\u00f74=56, 2\", \"596\u00f74=149, 0\"],\n  [\"241\u00f78=30, 1\", \"856\u00f75=171, 1\"],\n  [\"231\u00f74=57, 3\", \"105\u00f75=21, 0\"],\n  [\"216\u00f74=54, 0\", \"365\u00f79=40, 5\"],\n  [\"992\u00f76=165, 2\", \"541\u00f76=90, 1\"],\n  [\"467\u00f75=93, 2\", \"924\u00f78=115, 4\"],\n  [\"750\u00f75=150, 0\", \"464\u00f75=92, 4\"],\n  [\"663\u00f72=331, 1\", \"552\u00f78=69, 0\"],\n  [\"170\u00f78=21, 2\", \"368\u00f75=73, 3\"],\n  [\"419\u00f79=46, 5\", \"733\u00f74=183, 1\"],\n  [\"552\u00f75=110, 2\", \"566\u00f73=188, 2\"],\n  [\"130\u00f72=65, 0\", \"371\u00f79=41, 2\"],\n  [\"863\u00f75=172, 3\", \"803\u00f72=401, 1\"],\n  [\"301\u00f78=37, 5\", \"576\u00f72=288, 0\"],\n  [\"206\u00f79=22, 8\", \"407\u00f75=81, 2\"],\n  [\"564\u00f72=282, 0\", \"295\u00f72=147, 1\"],\n  [\"711\u00f72=355, 1\", \"791\u00f74=197, 3\"],\n  [\"483\u00f74=120, 3\", \"712\u00f72=356, 0\"],\n  [\"354\u00f76=59, 0\", \"784\u00f74=196, 0\"],\n  [\"467\u00f78=58, 3\", \"468\u00f78=58, 4\"],\n  [\"297\u00f74=74, 1\", \"795\u00f73=265, 0\"],\n  [\"652\u00f73=217, 1\", \"665\u00f73=221, 2\"],\n  [\"876\u00f72=438, 0\", \"606\u00f77=86, 4\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Could not find text: ${oldText}`);\n  }\n\n  results.items[0].insertText(newText, \"Replace\");\n}\n\nawait context.sync();\n", "ps1": "# Replace the 25 \"dividend\u00f7divisor=quotient, remainder\" answer strings in the\n# worksheet table with the new values from the commit.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"667\u00f76=111, 1\", \"878\u00f78=109, 6\"),\n    @(\"909\u00f72=454, 1\", \"827\u00f72=413, 1\"),\n    @(\"226\u00f74=56, 2\", \"596\u00f74=149, 0\"),\n    @(\"241\u00f78=30, 1\", \"856\u00f75=171, 1\"),\n    @(\"231\u00f74=57, 3\", \"105\u00f75=21, 0\"),\n    @(\"216\u00f74=54, 0\", \"365\u00f79=40, 5\"),\n    @(\"992\u00f76=165, 2\", \"541\u00f76=90, 1\"),\n    @(\"467\u00f75=93, 2\", \"924\u00f78=115, 4\"),\n    @(\"750\u00f75=150, 0\", \"464\u00f75=92, 4\"),\n    @(\"663\u00f72=331, 1\", \"552\u00f78=69, 0\"),\n    @(\"170\u00f78=21, 2\", \"368\u00f75=73, 3\"),\n    @(\"419\u00f79=46, 5\", \"733\u00f74=183, 1\"),\n    @(\"552\u00f75=110, 2\", \"566\u00f73=188, 2\"),\n    @(\"130\u00f72=65, 0\", \"371\u00f79=41, 2\"),\n    @(\"863\u00f75=172, 3\", \"803\u00f72=401, 1\"),\n    @(\"301\u00f78=37, 5\", \"576\u00f72=288, 0\"),\n    @(\"206\u00f79=22, 8\", \"407\u00f75=81, 2\"),\n    @(\"564\u00f72=282, 0\", \"295\u00f72=147, 1\"),\n    @(\"711\u00f72=355, 1\", \"791\u00f74=197, 3\"),\n    @(\"483\u00f74=120, 3\", \"712\u00f72=356, 0\"),\n    @(\"354\u00f76=59, 0\", \"784\u00f74=196, 0\"),\n    @(\"467\u00f78=58, 3\", \"468\u00f78=58, 4\"),\n    @(\"297\u00f74=74, 1\", \"795\u00f73=265, 0\"),\n    @(\"652\u00f73=217, 1\", \"665\u00f73=221, 2\"),\n    @(\"876\u00f72=438, 0\", \"606\u00f77=86, 4\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $true, $find.Replacement.Text, 2)\n}\n"}
